$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 1037
$ws.Range("F6").Value = 372
$ws.Range("F9").Value = 1487
$ws.Range("F11").Value = 1389
$ws.Range("F12").Value = 3035
$ws.Range("F13").Value = 521
$ws.Range("F14").Value = 1690
$ws.Range("F15").Value = 1374
$ws.Range("F17").Value = 252
$ws.Range("F18").Value = 1426
$ws.Range("F19").Value = 272
$ws.Range("F21").Value = 1158
$ws.Range("F22").Value = 31
$ws.Range("F23").Value = 417
$ws.Range("F24").Value = 36
$ws.Range("F25").Value = 3595
$ws.Range("F28").Value = 1588
$ws.Range("F29").Value = 39

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 41
$ws.Range("F5").Value = 16
$ws.Range("F6").Value = 57
$ws.Range("F7").Value = 8
$ws.Range("F8").Value = 33
$ws.Range("F9").Value = 27

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 22

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 22
$ws.Range("F8").Value = 41
$ws.Range("F9").Value = 16
$ws.Range("F10").Value = 57
$ws.Range("F12").Value = 8
$ws.Range("F13").Value = 33
$ws.Range("F14").Value = 27
$ws.Range("F16").Value = 1037
$ws.Range("F17").Value = 372
$ws.Range("F20").Value = 1487
$ws.Range("F22").Value = 1389
$ws.Range("F23").Value = 3035
$ws.Range("F24").Value = 521
$ws.Range("F25").Value = 1690
$ws.Range("F26").Value = 1374
$ws.Range("F28").Value = 252
$ws.Range("F29").Value = 1426
$ws.Range("F30").Value = 272
$ws.Range("F34").Value = 1158
$ws.Range("F35").Value = 31
$ws.Range("F36").Value = 417
$ws.Range("F37").Value = 36
$ws.Range("F38").Value = 3595
$ws.Range("F41").Value = 1588
$ws.Range("F44").Value = 39
